# Natmi following Dr Hou advice
# Adds a new cell-type cluster ("sCs") to the Col9a2 -> Mag ligand-receptor
# table, expanding the original single-row result (FAPs -> Col9a2/Mag -> M2)
# into the full 2x2 combination of sending/target clusters {FAPs, sCs}.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A..T correspond to the fixed header row already present:
#  A Sending cluster            K Receptor-expressing cells
#  B Ligand symbol              L Receptor detection rate
#  C Receptor symbol            M Receptor average expression value
#  D Target cluster             N Receptor total expression value
#  E Ligand-expressing cells    O Receptor derived specificity (avg)
#  F Ligand detection rate      P Receptor derived specificity (total)
#  G Ligand average expr value  Q Edge average expression weight
#  H Ligand total expr value    R Edge total expression weight
#  I Ligand derived spec (avg)  S Edge average expression derived specificity
#  J Ligand derived spec (tot)  T Edge total expression derived specificity

$rows = @(
    @("FAPs","Col9a2","Mag","M2",3,1,0.8244333333333334,2.4733,0.9386266517901193,0.9386266517901193,3,1,0.3133113333333333,0.939934,0.4010297802586483,0.4010297802586483,0.2583043069111111,2.3247387622,0.3764172399123023,0.3764172399123023),
    @("FAPs","Col9a2","Mag","sCs",3,1,0.8244333333333334,2.4733,0.9386266517901193,0.9386266517901193,2,0.6666666666666666,0.4679556666666667,1.403867,0.5989702197413518,0.5989702197413518,0.3857982501222222,3.4721842511,0.562209411877817,0.562209411877817),
    @("sCs","Col9a2","Mag","M2",1,0.3333333333333333,0.05390666666666667,0.16172,0.06137334820988076,0.06137334820988076,3,1,0.3133113333333333,0.939934,0.4010297802586483,0.4010297802586483,0.01688956960888889,0.15200612648,0.02461254034634599,0.02461254034634599),
    @("sCs","Col9a2","Mag","sCs",1,0.3333333333333333,0.05390666666666667,0.16172,0.06137334820988076,0.06137334820988076,2,0.6666666666666666,0.4679556666666667,1.403867,0.5989702197413518,0.5989702197413518,0.02522593013777778,0.22703337124,0.03676080786353478,0.03676080786353478)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    for ($j = 0; $j -lt $data.Count; $j++) {
        $col = $j + 1
        $ws.Cells.Item($r, $col).Value = $data[$j]
    }
}
